$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - exact text values, preserved via apostrophe-prefix
# to prevent Excel auto-converting numeric-looking strings to Number type,
# then Style reset to Normal to avoid leaving a quotePrefix style behind.
$dPrices = @{
    2 = "265.37"
    4 = "6.276"
    5 = "0.06153"
    6 = "3.575"
    7 = "6.662"
    8 = "1.343"
    9 = "0.8294"
    11 = "0.1593"
    12 = "0.08276"
    13 = "0.03424"
    14 = "0.03139"
    15 = "0.09256"
    16 = "3.896"
    17 = "0.001708"
    18 = "0.04891"
    19 = "0.006218"
    20 = "0.005275"
    21 = "0.001090"
    22 = "0.0001501"
    23 = "3.768"
    24 = "2.289"
    26 = "0.1238"
    27 = "0.0002681"
    40 = "0.04619"
    42 = "0.1137"
    43 = "0.003401"
    44 = "0.01081"
    45 = "0.00006140"
    47 = "0.7004"
    48 = "0.1937"
    49 = "0.00002101"
    50 = "0.01241"
}
foreach ($r in $dPrices.Keys) {
    $ws.Cells.Item($r, 4).Value = "'" + $dPrices[$r]
    $ws.Cells.Item($r, 4).Style = "Normal"
}

# Column G (Hora) updates - every data row 2-51 goes from "3" to "4"
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 7).Value = "'" + "4"
    $ws.Cells.Item($r, 7).Style = "Normal"
}
